$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy style from E1 (bold header style) and set text
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Timestamp values for rows 2-12 (stored as text, not as dates)
$timestamps = @(
    "2021-10-05 13:39:42.289565",
    "2021-10-05 13:39:42.289575",
    "2021-10-05 13:39:42.289579",
    "2021-10-05 13:39:42.289581",
    "2021-10-05 13:39:42.289584",
    "2021-10-05 13:39:42.289587",
    "2021-10-05 13:39:42.289590",
    "2021-10-05 13:39:42.289592",
    "2021-10-05 13:39:42.289595",
    "2021-10-05 13:39:42.289598",
    "2021-10-05 13:39:42.289601"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
